$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.276.91"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "1.773.63"
$ws.Range("E3").Value = "  +3.62%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9979"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5214"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +10.51%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3617"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.59"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07339"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.083"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9978"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.059"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.59%  "
$ws.Range("D15").Value = "1.770.25"
$ws.Range("E15").Value = "  +3.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.967"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001044"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06415"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9979"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.848"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.38%  "
$ws.Range("D23").Value = "27.349.08"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.062"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.33%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.336"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.49%  "
$ws.Range("B29").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C29").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D29").Value = "1.971.06"
$ws.Range("E29").Value = "  +3.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.062"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09749"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.554"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.597"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02226"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05978"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.19"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2029"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.826"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6116"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.59%  "
$ws.Range("E41").Value = "  +2.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.977"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.142"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5753"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.620"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.883"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.108"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06698"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.96%  "
